$wb = $excel.ActiveWorkbook

# =============================================================================
# "K10+" sheet — fix a mis-typed collection occasion (row 5) and add the
# missing collection occasion (row 6) that was dropped from the export.
# =============================================================================
$k10 = $wb.Worksheets.Item("K10+")

# Row 5 corrections (measure_date, reason_for_collection, k10p_item11-13)
$k10.Cells.Item(5, 4).Value = 20052016   # D5 measure_date
$k10.Cells.Item(5, 5).Value = 2          # E5 reason_for_collection
$k10.Cells.Item(5, 16).Value = 99        # P5 k10p_item11
$k10.Cells.Item(5, 17).Value = 99        # Q5 k10p_item12
$k10.Cells.Item(5, 18).Value = 99        # R5 k10p_item13

# New row 6 - additional K10+ collection occasion
$k10.Cells.Item(6, 1).Value = "PHN999:NFP01"  # A6 Version
$k10.Cells.Item(6, 2).Value = "CO04"          # B6 collection_occasion_key
$k10.Cells.Item(6, 3).Value = "E01"           # C6 episode_key
$k10.Cells.Item(6, 4).Value = 18062016        # D6 measure_date
$k10.Cells.Item(6, 5).Value = 3               # E6 reason_for_collection
$k10.Cells.Item(6, 6).Value = 1                # F6 k10p_item1
$k10.Cells.Item(6, 7).Value = 1                # G6 k10p_item2
$k10.Cells.Item(6, 8).Value = 1                # H6 k10p_item3
$k10.Cells.Item(6, 9).Value = 1                # I6 k10p_item4
$k10.Cells.Item(6, 10).Value = 1               # J6 k10p_item5
$k10.Cells.Item(6, 11).Value = 1               # K6 k10p_item6
$k10.Cells.Item(6, 12).Value = 1               # L6 k10p_item7
$k10.Cells.Item(6, 13).Value = 1               # M6 k10p_item8
$k10.Cells.Item(6, 14).Value = 1               # N6 k10p_item9
$k10.Cells.Item(6, 15).Value = 1               # O6 k10p_item10
$k10.Cells.Item(6, 16).Value = 1               # P6 k10p_item11
$k10.Cells.Item(6, 17).Value = 3               # Q6 k10p_item12
$k10.Cells.Item(6, 18).Value = 1               # R6 k10p_item13
$k10.Cells.Item(6, 19).Value = 9               # S6 k10p_item14
$k10.Cells.Item(6, 20).Value = 99              # T6 k10p_score
$k10.Cells.Item(6, 21).Value = "tag1"          # U6 k10p_tags

# P6:R6 carry an explicit black font colour in the source workbook
$k10.Range("P6:R6").Font.Color = 0

# =============================================================================
# "SDQ" sheet — correct the SDQ-P item codes on row 3 (9 = not-answered
# placeholder was wrong; the real value recorded on the paper form was 8).
# =============================================================================
$sdq = $wb.Worksheets.Item("SDQ")
$sdq.Cells.Item(3, 40).Value = 8   # AN3
$sdq.Cells.Item(3, 41).Value = 8   # AO3
$sdq.Cells.Item(3, 45).Value = 8   # AS3
$sdq.Cells.Item(3, 46).Value = 8   # AT3
$sdq.Cells.Item(3, 47).Value = 8   # AU3
$sdq.Cells.Item(3, 48).Value = 8   # AV3

# =============================================================================
# "Practitioners" sheet — correct row 3 and fill in the previously blank
# atsi_cultural_training / practitioner_active cells.
# =============================================================================
$prac = $wb.Worksheets.Item("Practitioners")
$prac.Cells.Item(3, 4).Value = 3   # D3 atsi_cultural_training
$prac.Cells.Item(4, 4).Value = 1   # D4 atsi_cultural_training
$prac.Cells.Item(5, 8).Value = 0   # H5 practitioner_active

# =============================================================================
# Restore each sheet's on-screen selection, matching where the author was
# last working, then leave "Practitioners" as the active tab.
# =============================================================================
$sc = $wb.Worksheets.Item("Service Contacts")
$sc.Activate() | Out-Null
$sc.Range("E13").Select() | Out-Null

$k10.Activate() | Out-Null
$k10.Range("L18").Select() | Out-Null

$sdq.Activate() | Out-Null
$sdq.Range("H14").Select() | Out-Null

$prac.Activate() | Out-Null
$prac.Range("H3").Select() | Out-Null
